# Add documentation examples (E0044-E0049) to the "Example" sheet.
# Column layout: A=Term, B=Title, C=Description, D=Source, E=SourceFormat,
#                F=SourceType, G=Concepts, H=Reference, I=Status, J=Date,
#                K=Contributor
#
# New rows reuse the same formatting (cell styles, SourceFormat/SourceType/
# Status/Date/Contributor) as the last existing data row (44), so that row
# is copied as a template for each new row before the per-row text values
# are overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

$rows = @(
    @{
        A = "E0044"
        B = "Specifying personal data"
        C = "This example shows how personal data can be specified by for a work email address in four methods: (1) by using the PD taxonomy, (2) extending a concept from PD taxonomy,(3) specifying the exact value of email address, and (4) using external vocabularies such as schema.org. In all of these, it is strongly recommended to include a DPV concept (i.e. declare it as Personal Data) for interoperability, and to align it with the DPV and PD taxonomies where relevant (e.g. to mention it is an email address, or to specify it is identifier, or that it is sensitive - which it is not in this case)"
        D = "E0044.ttl"
        G = "dpv:PersonalData,dpv:hasPersonalData"
    },
    @{
        A = "E0045"
        B = "Indicating data belongs to sensitive or special category"
        C = "This example shows how data can be indicated as being sensitive or belonging to special category. It also shows the use of PD extension which provides a taxonomy of special categories."
        D = "E0045.ttl"
        G = "dpv:SensitivePersonalData,dpv:SpecialCategoryPersonalData"
    },
    @{
        A = "E0046"
        B = "Indicating data being collected and derived"
        C = "This example shows a process which first collects email address from the data subject, and then uses it to derive an account identifier. The seeming duplication in information across processing, personal data category, and data source actually represents three distinct concepts - which can be used in various ways for data governance, or legal compliance e.g. to retrieve all data which is collected or to ensure all collected data has a source."
        D = "E0046.ttl"
        G = "dpv:CollectedPersonalData,dpv:ProvidedPersonalData,dpv:DerivedPersonalData"
    },
    @{
        A = "E0047"
        B = "Indicating processing conditions for duration and location"
        C = "This example shows processing conditions where the use (of data or technology) takes place over 6 months and in Ireland (IE) and Netherlands (NL)"
        D = "E0047.ttl"
        G = "dpv:ProcessingCondition,dpv:ProcessingDuration,dpv:ProcessingLocation"
    },
    @{
        A = "E0048"
        B = "Indicating storage conditions for duration, location, deletion, and restoration"
        C = "This example shows storage conditions for a 'store' processing operation. It has a duration valid until the event 'account closure' occurs. It has a storage location situated in Ireland (IE) and Netherlands (NL). The deletion occurs 6 months after the event 'account closure'. And restoration is implemented by using (stored) data located in backup systems in Ireland (IE)."
        D = "E0048.ttl"
        G = "dpv:StorageCondition,dpv:StorageDuration,dpv:StorageLocation,dpv:StorageDeletion,dpv:StorageRestoration"
    },
    @{
        A = "E0049"
        B = "Indicating data volume, geo-location coverage, data subject scale, and a processing scale"
        C = "This example shows how data volume, data subject scale, and geographic scale can be indicated along with optional information about the exact values involved. It also shows how a qualified scale can be provided as 'processing scale' based on some criteria (not described here)."
        D = "E0049.ttl"
        G = "dpv:ProcessingScale,dpv:GeographicCoverage,dpv:DataVolume,dpv:DataSubjectScale"
    }
)

$templateRow = 44
$startRow = 45

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Copy formatting (and placeholder values) from the template row, split
    # around column H so no blank H cell is materialised (rows 31-44 in the
    # original sheet likewise omit H entirely).
    $ws.Range("A$templateRow`:G$templateRow").Copy($ws.Range("A$r`:G$r"))
    $ws.Range("I$templateRow`:K$templateRow").Copy($ws.Range("I$r`:K$r"))

    $ws.Cells.Item($r, 1).Value = $data.A   # Term
    $ws.Cells.Item($r, 2).Value = $data.B   # Title
    $ws.Cells.Item($r, 3).Value = $data.C   # Description
    $ws.Cells.Item($r, 4).Value = $data.D   # Source
    $ws.Cells.Item($r, 7).Value = $data.G   # Concepts
}
